# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 1174
    $ws.Range("F6").Value = 169
    $ws.Range("F10").Value = 5460
    $ws.Range("F11").Value = 4871
}
